$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.691.90'
$ws.Range('E2').Value = '  +3.06%  '
$ws.Range('D3').Value = '2.446.55'
$ws.Range('E3').Value = '  +1.90%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'578.11"
$ws.Range('E5').Value = '  +2.85%  '
$ws.Range('D6').Value = "'145.81"
$ws.Range('E6').Value = '  +2.99%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').Value = '2.445.85'
$ws.Range('E10').Value = '  +2.69%  '
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').Value = "'5.24"
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').Value = "'0.353"
$ws.Range('E13').Value = '  +2.86%  '
$ws.Range('D14').Value = "'28.29"
$ws.Range('E14').Value = '  +8.18%  '
$ws.Range('E15').Value = '  +5.60%  '
$ws.Range('D17').Value = '62.538.88'
$ws.Range('E17').Value = '  +3.03%  '
$ws.Range('D18').Value = '2.459.74'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').Value = "'7.85"
$ws.Range('E19').Value = '  -3.95%  '
$ws.Range('D20').Value = "'10.94"
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('D21').Value = "'328.25"
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('E23').Value = '  +10.40%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = "'65.62"
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').Value = "'639.86"
$ws.Range('E26').Value = '  +11.97%  '
$ws.Range('E27').Value = '  +16.95%  '
$ws.Range('D28').Value = "'8.47"
$ws.Range('E28').Value = '  +5.26%  '
$ws.Range('D29').Value = '0.0₃0983'
$ws.Range('E29').Value = '  +4.68%  '
$ws.Range('E30').Value = '  +85.46%  '
$ws.Range('D31').Value = '2.575.16'
$ws.Range('D32').Value = "'8.20"
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('D33').Value = "'1.45"
$ws.Range('E33').Value = '  +8.35%  '
$ws.Range('D34').Value = "'1.88"
$ws.Range('E34').Value = '  +3.90%  '
$ws.Range('D35').Value = "'0.138"
$ws.Range('E35').Value = '  +4.19%  '
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E38').Value = '  +3.47%  '
$ws.Range('E39').Value = '  +6.10%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = "'0.374"
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = "'152.80"
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = "'18.67"
$ws.Range('E42').Value = '  +1.99%  '
$ws.Range('E43').Value = '  +6.28%  '
$ws.Range('E44').Value = '  +5.17%  '
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('E47').Value = '  +27.77%  '
$ws.Range('D48').Value = "'145.43"
$ws.Range('E48').Value = '  +2.34%  '
$ws.Range('D49').Value = "'3.60"
$ws.Range('E49').Value = '  +2.76%  '
$ws.Range('D50').Value = "'20.64"
$ws.Range('E50').Value = '  +6.81%  '
$ws.Range('E51').Value = '  +3.05%  '
